# CAN_scaling_mapping.xlsx edit: remove shipping / aviation related rows
# per commit message "Remove shipping and aviation, other changes."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")

# Row 36 (Air Transportation / 1A3ai_International-aviation):
#   clear the "Air" scaling_sector (B36) and the Notes cell (D36)
$ws.Range("B36").ClearContents()
$ws.Range("D36").ClearContents()

# Row 37 (1A3aii_Domestic-aviation):
#   clear the "Air" scaling_sector (B37)
$ws.Range("B37").ClearContents()

# Row 44 (Marine Transportation):
#   clear the "Marine" scaling_sector (B44) and the
#   "1A3dii_Domestic-naviation" ceds_sector (C44)
$ws.Range("B44").ClearContents()
$ws.Range("C44").ClearContents()

# Update the sheet view to reflect the scroll/selection state seen after
# the edit session.
$ws.Activate()
$ws.Range("D37").Select()
